$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.15%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.21%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.304"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.90%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'4.32%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.588"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.91%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.340"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'27.78%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.646"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.59%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1280"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.09%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1974"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'6.15%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09690"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.02%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04699"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.26%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-0.19%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001320"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.04%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04202"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.75%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005922"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.79%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.349"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.442"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'4.80%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3520"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'5.04%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.117"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.08%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1382"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.14%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3081"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-3.04%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001295"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.26%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004321"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.99%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.04%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003533"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02722"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'8.37%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'12.83%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01077"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'85.33%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008017"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.55%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1469"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.06%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007898"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.37%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007873"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.05%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'16.62%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007072"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'5.86%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05577"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'27.26%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003994"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.93%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002096"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.05%"
$ws.Range("E51").Style = "Normal"
